# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Apio" (row 120 and 121), shifting the
# previously existing rows 120-132 down to 122-134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 120:132 down by two rows (to 122:134), duplicating the
# formatting (e.g. the date number format on column D) into the freshly
# inserted blank rows.
$ws.Rows("120:121").Insert()

# --- New row 120 -----------------------------------------------------------
$ws.Range("A120").Value = 9
$ws.Range("B120").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C120").Value = "Metropolitana"
$ws.Range("D120").Value = "10/05/2021"
$ws.Range("E120").Value = 13
$ws.Range("F120").Value = 100112017
$ws.Range("G120").Value = "Apio"
$ws.Range("H120").Value = "Americana (o)"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 52
$ws.Range("K120").Value = 8000
$ws.Range("L120").Value = 9000
$ws.Range("M120").Value = 8500
$ws.Range("N120").Value = "`$/docena de matas"
$ws.Range("O120").Value = "Región de Coquimbo"
$ws.Range("P120").Value = 1417
$ws.Range("Q120").Value = 6
$ws.Range("R120").Value = "Hortaliza"

# --- New row 121 -----------------------------------------------------------
$ws.Range("A121").Value = 9
$ws.Range("B121").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C121").Value = "Metropolitana"
$ws.Range("D121").Value = "10/05/2021"
$ws.Range("E121").Value = 13
$ws.Range("F121").Value = 100112017
$ws.Range("G121").Value = "Apio"
$ws.Range("H121").Value = "Americana (o)"
$ws.Range("I121").Value = "Segunda"
$ws.Range("J121").Value = 34
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 7000
$ws.Range("M121").Value = 6500
$ws.Range("N121").Value = "`$/docena de matas"
$ws.Range("O121").Value = "Región de Coquimbo"
$ws.Range("P121").Value = 1083
$ws.Range("Q121").Value = 6
$ws.Range("R121").Value = "Hortaliza"
